{"js": "// Fix the lab exercise date in the title block: \"Lab Exercise 11/20/2022\"\n// (spread across several runs as \"Lab Exercise 1\" + \"1/20\" + \"202\" + \"2\")\n// should read \"Lab Exercise 11/30/2022\" (the \"1/20\" run becomes \"1/30/\").\nconst body = context.document.body;\n\nconst results = body.search(\"1/20\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the run's text in place so its bold / 28pt formatting is kept.\n  results.items[0].insertText(\"1/30/\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Fix the lab exercise date in the title block: \"Lab Exercise 11/20/2022\"\n# (spread across several runs as \"Lab Exercise 1\" + \"1/20\" + \"202\" + \"2\")\n# should read \"Lab Exercise 11/30/2022\" (the \"1/20\" run becomes \"1/30/\").\n$d = $word.ActiveDocument\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.MatchCase = $true\n$r.Find.MatchWildcards = $false\n\nif ($r.Find.Execute(\"1/20\")) {\n    # Replace just the matched range's text so its bold / 28pt run\n    # formatting (inherited from the original run) is preserved.\n    $r.Text = \"1/30/\"\n}\n"}
